$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 49, shifting existing rows 49-60 down to 50-61
$ws.Rows.Item(49).Insert()

# Populate the new row 49 with the new weekly data point
$ws.Cells.Item(49, 1).Value = 11
$ws.Cells.Item(49, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(49, 3).Value = "Bíobío"
$ws.Cells.Item(49, 4).Value = 45135
$ws.Cells.Item(49, 4).NumberFormat = $ws.Cells.Item(50, 4).NumberFormat
$ws.Cells.Item(49, 5).Value = 8
$ws.Cells.Item(49, 6).Value = 100112026
$ws.Cells.Item(49, 7).Value = "Haba"
$ws.Cells.Item(49, 8).Value = "Sin especificar"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 50
$ws.Cells.Item(49, 11).Value = 16000
$ws.Cells.Item(49, 12).Value = 16000
$ws.Cells.Item(49, 13).Value = 16000
$ws.Cells.Item(49, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(49, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(49, 16).Value = 640
$ws.Cells.Item(49, 17).Value = 25
$ws.Cells.Item(49, 18).Value = "Hortaliza"
